$d = $word.ActiveDocument

$replacements = @(
    @("69×66=", "89×34="),
    @("86×44=", "96×78="),
    @("34×47=", "42×40="),
    @("29×56=", "41×80="),
    @("85×41=", "75×41="),
    @("79×99=", "44×80="),
    @("45×65=", "42×12="),
    @("86×28=", "12×38="),
    @("56×51=", "82×69="),
    @("59×83=", "76×46="),
    @("98×43=", "84×51="),
    @("65×28=", "60×49="),
    @("83×34=", "71×74="),
    @("34×23=", "28×66="),
    @("76×57=", "32×29="),
    @("37×69=", "17×49="),
    @("97×93=", "45×75="),
    @("33×33=", "38×45="),
    @("43×95=", "63×42="),
    @("44×33=", "29×81="),
    @("55×20=", "53×99="),
    @("61×94=", "11×79="),
    @("69×95=", "21×28="),
    @("11×23=", "59×44="),
    @("20×86=", "60×85=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
